$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-22 05:08:41"
$wsZh.Range("H2").Value = "2016-03-22 05:09:03"

# de-de sheet: update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-22 05:08:44"
$wsDe.Range("H2").Value = "2016-03-22 05:09:09"
